$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "fLD:IWoWRJ[ru3R9I_P}r"
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = "ir"

$ws.Range("C5").Select()
